function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$wsRef = $wb.Worksheets.Item(1)   # style-copy source (header/index cell formatting)

# ---------------------------------------------------------------
# Step 1: duplicate the existing "总计" sheet (6th sheet) so the copy
#         keeps its quarterly summary data + page/outline formatting,
#         then rename the original to "2022-Q1" and the copy back to
#         "总计" (this also reproduces the sheetId numbering: the
#         renamed original keeps sheetId 6, the new copy becomes 7).
# ---------------------------------------------------------------
$sheetTotalOrig = $wb.Worksheets.Item(6)
$sheetTotalOrig.Copy([System.Reflection.Missing]::Value, $sheetTotalOrig)
$sheetTotalOrig.Name = "2022-Q1"
$sheetTotalNew = $wb.Worksheets.Item(7)
$sheetTotalNew.Name = "总计"

# ---- Rebuild the "2022-Q1" sheet from scratch with fund holding detail ----
$wsQ1 = $sheetTotalOrig
$wsQ1.Cells.Clear()

Set-TextCell $wsQ1 "B1" "基金代码"
Set-TextCell $wsQ1 "C1" "基金名称"
Set-TextCell $wsQ1 "D1" "基金规模"
Set-TextCell $wsQ1 "E1" "股票总仓位"
Set-TextCell $wsQ1 "F1" "仓位占比"
Set-TextCell $wsQ1 "G1" "持有市值(亿元)"
Set-TextCell $wsQ1 "H1" "仓位排名"

# Header style (bold + border + centered), matching the other quarter sheets
$wsRef.Range("B1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ1.Range("A2").Value = 0
Set-TextCell $wsQ1 "B2" "162607"
Set-TextCell $wsQ1 "C2" "景顺长城资源垄断混合(LOF)"
Set-TextCell $wsQ1 "D2" "17.99"
Set-TextCell $wsQ1 "E2" "93.81"
Set-TextCell $wsQ1 "F2" "8.32"
Set-TextCell $wsQ1 "G2" "1.4968"
$wsQ1.Range("H2").Value = 6

$wsQ1.Range("A3").Value = 1
Set-TextCell $wsQ1 "B3" "000772"
Set-TextCell $wsQ1 "C3" "景顺长城中国回报灵活配置混合"
Set-TextCell $wsQ1 "D3" "16.56"
Set-TextCell $wsQ1 "E3" "93.73"
Set-TextCell $wsQ1 "F3" "7.90"
Set-TextCell $wsQ1 "G3" "1.3082"
$wsQ1.Range("H3").Value = 6

$wsQ1.Range("A4").Value = 2
Set-TextCell $wsQ1 "B4" "010201"
Set-TextCell $wsQ1 "C4" "农银汇理智增一年定期开放混合"
Set-TextCell $wsQ1 "D4" "9.82"
Set-TextCell $wsQ1 "E4" "62.05"
Set-TextCell $wsQ1 "F4" "8.11"
Set-TextCell $wsQ1 "G4" "0.7964"
$wsQ1.Range("H4").Value = 3

$wsQ1.Range("A5").Value = 3
Set-TextCell $wsQ1 "B5" "001656"
Set-TextCell $wsQ1 "C5" "农银汇理中国优势灵活配置混合"
Set-TextCell $wsQ1 "D5" "2.92"
Set-TextCell $wsQ1 "E5" "68.89"
Set-TextCell $wsQ1 "F5" "8.66"
Set-TextCell $wsQ1 "G5" "0.2529"
$wsQ1.Range("H5").Value = 1

$wsQ1.Range("A6").Value = 4
Set-TextCell $wsQ1 "B6" "660005"
Set-TextCell $wsQ1 "C6" "农银中小盘混合"
Set-TextCell $wsQ1 "D6" "7.31"
Set-TextCell $wsQ1 "E6" "68.62"
Set-TextCell $wsQ1 "F6" "2.17"
Set-TextCell $wsQ1 "G6" "0.1586"
$wsQ1.Range("H6").Value = 8

$wsQ1.Range("A7").Value = 5
Set-TextCell $wsQ1 "B7" "121006"
Set-TextCell $wsQ1 "C7" "国投瑞银稳健增长混合"
Set-TextCell $wsQ1 "D7" "7.34"
Set-TextCell $wsQ1 "E7" "65.10"
Set-TextCell $wsQ1 "F7" "2.06"
Set-TextCell $wsQ1 "G7" "0.1512"
$wsQ1.Range("H7").Value = 10

$wsQ1.Range("A8").Value = 6
Set-TextCell $wsQ1 "B8" "011349"
Set-TextCell $wsQ1 "C8" "淳厚现代服务业股票A"
Set-TextCell $wsQ1 "D8" "3.58"
Set-TextCell $wsQ1 "E8" "81.51"
Set-TextCell $wsQ1 "F8" "2.96"
Set-TextCell $wsQ1 "G8" "0.1060"
$wsQ1.Range("H8").Value = 5

$wsQ1.Range("A9").Value = 7
Set-TextCell $wsQ1 "B9" "007811"
Set-TextCell $wsQ1 "C9" "淳厚信泽灵活配置混合A"
Set-TextCell $wsQ1 "D9" "4.32"
Set-TextCell $wsQ1 "E9" "74.11"
Set-TextCell $wsQ1 "F9" "2.25"
Set-TextCell $wsQ1 "G9" "0.0972"
$wsQ1.Range("H9").Value = 9

$wsQ1.Range("A10").Value = 8
Set-TextCell $wsQ1 "B10" "020015"
Set-TextCell $wsQ1 "C10" "国泰区位优势混合"
Set-TextCell $wsQ1 "D10" "1.44"
Set-TextCell $wsQ1 "E10" "80.83"
Set-TextCell $wsQ1 "F10" "5.49"
Set-TextCell $wsQ1 "G10" "0.0791"
$wsQ1.Range("H10").Value = 2

$wsQ1.Range("A11").Value = 9
Set-TextCell $wsQ1 "B11" "004099"
Set-TextCell $wsQ1 "C11" "前海开源沪港深景气行业精选灵活配置混合"
Set-TextCell $wsQ1 "D11" "0.41"
Set-TextCell $wsQ1 "E11" "93.07"
Set-TextCell $wsQ1 "F11" "9.37"
Set-TextCell $wsQ1 "G11" "0.0384"
$wsQ1.Range("H11").Value = 3

$wsQ1.Range("A12").Value = 10
Set-TextCell $wsQ1 "B12" "007812"
Set-TextCell $wsQ1 "C12" "淳厚信泽灵活配置混合C"
Set-TextCell $wsQ1 "D12" "1.15"
Set-TextCell $wsQ1 "E12" "74.11"
Set-TextCell $wsQ1 "F12" "2.25"
Set-TextCell $wsQ1 "G12" "0.0259"
$wsQ1.Range("H12").Value = 9

$wsQ1.Range("A13").Value = 11
Set-TextCell $wsQ1 "B13" "011824"
Set-TextCell $wsQ1 "C13" "浙商汇金量化臻选股票型证券投资基金A"
Set-TextCell $wsQ1 "D13" "1.54"
Set-TextCell $wsQ1 "E13" "92.80"
Set-TextCell $wsQ1 "F13" "1.53"
Set-TextCell $wsQ1 "G13" "0.0236"
$wsQ1.Range("H13").Value = 4

$wsQ1.Range("A14").Value = 12
Set-TextCell $wsQ1 "B14" "011350"
Set-TextCell $wsQ1 "C14" "淳厚现代服务业股票C"
Set-TextCell $wsQ1 "D14" "0.63"
Set-TextCell $wsQ1 "E14" "81.51"
Set-TextCell $wsQ1 "F14" "2.96"
Set-TextCell $wsQ1 "G14" "0.0186"
$wsQ1.Range("H14").Value = 5

$wsQ1.Range("A15").Value = 13
Set-TextCell $wsQ1 "B15" "002182"
Set-TextCell $wsQ1 "C15" "东兴蓝海财富灵活配置混合"
Set-TextCell $wsQ1 "D15" "0.24"
Set-TextCell $wsQ1 "E15" "86.93"
Set-TextCell $wsQ1 "F15" "4.39"
Set-TextCell $wsQ1 "G15" "0.0105"
$wsQ1.Range("H15").Value = 8

$wsQ1.Range("A16").Value = 14
Set-TextCell $wsQ1 "B16" "011825"
Set-TextCell $wsQ1 "C16" "浙商汇金量化臻选股票型证券投资基金C"
Set-TextCell $wsQ1 "D16" "0.47"
Set-TextCell $wsQ1 "E16" "92.80"
Set-TextCell $wsQ1 "F16" "1.53"
Set-TextCell $wsQ1 "G16" "0.0072"
$wsQ1.Range("H16").Value = 4

$wsQ1.Range("A17").Value = 15
Set-TextCell $wsQ1 "B17" "003717"
Set-TextCell $wsQ1 "C17" "中银量化精选灵活配置混合A"
Set-TextCell $wsQ1 "D17" "0.49"
Set-TextCell $wsQ1 "E17" "90.38"
Set-TextCell $wsQ1 "F17" "1.24"
Set-TextCell $wsQ1 "G17" "0.0061"
$wsQ1.Range("H17").Value = 5

$wsQ1.Range("A18").Value = 16
Set-TextCell $wsQ1 "B18" "001849"
Set-TextCell $wsQ1 "C18" "前海开源强势共识100强等权重股票"
Set-TextCell $wsQ1 "D18" "0.12"
Set-TextCell $wsQ1 "E18" "92.23"
Set-TextCell $wsQ1 "F18" "1.04"
Set-TextCell $wsQ1 "G18" "0.0012"
$wsQ1.Range("H18").Value = 9

$wsQ1.Range("A19").Value = 17
Set-TextCell $wsQ1 "B19" "010484"
Set-TextCell $wsQ1 "C19" "中银量化精选灵活配置混合C"
Set-TextCell $wsQ1 "D19" "0.01"
Set-TextCell $wsQ1 "E19" "90.38"
Set-TextCell $wsQ1 "F19" "1.24"
Set-TextCell $wsQ1 "G19" "0.0001"
$wsQ1.Range("H19").Value = 5

# Style column-A index cells (s="2", border/bold/center) like the other sheets
$wsRef.Range("A2").Copy()
$wsQ1.Range("A2:A19").PasteSpecial(-4122)

# ---- Update the "总计" sheet: insert the new 2022-Q1 row at the top ----
$wsTotal = $sheetTotalNew
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").Style = "Normal"
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 18
$wsTotal.Range("D2").Value = 4.58

# Restore the column-A running index (0,1,2,...) across all rows
for ($i = 0; $i -le 5; $i++) {
    $wsTotal.Range("A" + (2 + $i)).Value = $i
}
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

# Restore original active sheet/selection (first sheet), matching the source workbook
$wsRef.Activate()
$null = $wsRef.Range("A1").Select()

Write-Host "Edit complete"